$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.118.48"
$ws.Range("D3").Value = "1.832.13"
$ws.Range("E3").Value = "  -0.13%  "
$origStyle = $ws.Range("D4").Style
$ws.Range("D4").Value = "'0.9998"
$ws.Range("D4").Style = $origStyle
$ws.Range("E4").Value = "  +0.05%  "
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").Value = "'243.29"
$ws.Range("D5").Style = $origStyle
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").Value = "'0.6271"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("E7").Value = "  +0.05%  "
$origStyle = $ws.Range("D8").Style
$ws.Range("D8").Value = "'0.07502"
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = "  -0.92%  "
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").Value = "'0.2932"
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = "  +0.23%  "
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").Value = "'23.28"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  +2.95%  "
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").Value = "'0.07706"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  -0.48%  "
$ws.Range("D12").Value = "1.817.11"
$ws.Range("E12").Value = "  -0.80%  "
$origStyle = $ws.Range("D13").Style
$ws.Range("D13").Value = "'5.029"
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = "  +1.23%  "
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").Value = "'0.6692"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = "  +0.59%  "
$origStyle = $ws.Range("D15").Style
$ws.Range("D15").Value = "'82.89"
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = "  -0.27%  "
$origStyle = $ws.Range("D16").Style
$ws.Range("D16").Value = "'0.000009370"
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = "  -6.93%  "
$origStyle = $ws.Range("D17").Style
$ws.Range("D17").Value = "'5.999"
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = "  -1.17%  "
$ws.Range("D18").Value = "29.110.03"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").Value = "2.076.00"
$ws.Range("E19").Value = "  -0.39%  "
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").Value = "'223.32"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  -1.66%  "
$ws.Range("E22").Value = "  +0.03%  "
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").Value = "'7.145"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  -0.97%  "
$ws.Range("E24").Value = "  +0.03%  "
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").Value = "'160.32"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  +0.26%  "
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").Value = "'0.1401"
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = "  +1.03%  "
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").Value = "'8.514"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  +0.00%  "
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").Value = "'17.91"
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = "  -0.18%  "
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").Value = "'1.491"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  -0.39%  "
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").Value = "'0.05837"
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = "  +11.19%  "
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").Value = "'4.166"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  +1.49%  "
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").Value = "'4.130"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = "  +2.85%  "
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").Value = "'1.210"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  +1.33%  "
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").Value = "'0.7428"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  +0.84%  "
$ws.Range("E35").Value = "  -0.56%  "
$origStyle = $ws.Range("D36").Style
$ws.Range("D36").Value = "'1.140"
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = "  +0.21%  "
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").Value = "'2.671"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  -0.56%  "
$ws.Range("D38").Value = "1.233.67"
$ws.Range("E38").Value = "  -0.74%  "
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").Value = "'2.762"
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("E40").Value = "  -0.28%  "
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").Value = "'6.505"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  +2.30%  "
$ws.Range("E42").Value = "  -0.78%  "
$ws.Range("E43").Value = "  -0.02%  "
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").Value = "'102.24"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("D45").Value = "1.976.58"
$ws.Range("E45").Value = "  -0.29%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").Value = "'0.00000000125"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  +1.98%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").Value = "'66.04"
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = "  +2.78%  "
$ws.Range("E48").Value = "  -0.47%  "
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").Value = "'0.4075"
$ws.Range("D49").Style = $origStyle
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").Value = "'0.07507"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  +12.47%  "
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").Value = "'8.999"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  +1.23%  "
